# edit.ps1
# Applies the scraped-data update to the "serbia_prva-liga_2023-2024" sheet:
#  1. A number of existing rows (identified by match/date groupings) get their
#     match-data block (columns F:V) rotated among themselves - this reflects
#     the source site re-ordering matches that share the same date/time while
#     the underlying row (A:E - Indice/pais/torneio/temporada/data_partida)
#     keeps its original position.
#  2. Six brand-new match rows are appended at the end of the table (rows
#     141-146), extending the sheet's used range to A1:V146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-Block($ws, $row) {
    return $ws.Range("F$row`:V$row").Value2
}

function Set-Block($ws, $row, $data) {
    $ws.Range("F$row`:V$row").Value = $data
}

# Rotate the F:V content across the given ordered list of rows:
# new-content(rows[i]) = old-content(rows[i+1]) ; last row wraps to the first.
function Rotate-Cycle($ws, [int[]]$rows) {
    $blocks = @()
    foreach ($r in $rows) {
        $blocks += ,(Get-Block $ws $r)
    }
    $n = $rows.Count
    for ($i = 0; $i -lt $n; $i++) {
        $src = $blocks[($i + 1) % $n]
        Set-Block $ws $rows[$i] $src
    }
}

Rotate-Cycle $ws @(22, 23)
Rotate-Cycle $ws @(65, 67, 66)
Rotate-Cycle $ws @(73, 75)
Rotate-Cycle $ws @(84, 86, 85)
Rotate-Cycle $ws @(87, 89, 88)
Rotate-Cycle $ws @(103, 105, 104)
Rotate-Cycle $ws @(106, 107)
Rotate-Cycle $ws @(108, 109, 110)
Rotate-Cycle $ws @(114, 116, 115)
Rotate-Cycle $ws @(119, 120, 121, 122)
Rotate-Cycle $ws @(126, 130, 127, 129)

# Append the six new rows (141-146) at the bottom of the table.
# Columns A (Indice) and E (data_partida) reuse the formatting of the last
# existing data row (140) so the bold/border/centering and the datetime
# number-format stay consistent with the rest of the column.
$ws.Range("A140").Copy($ws.Range("A141:A146"))
$ws.Range("E140").Copy($ws.Range("E141:E146"))

$newRows = @(
    @{ Row = 141; A = 140; E = 45262.54166666666; F = "Metalac";          G = 2; H = "Kolubara";               I = 2;
       J = 2.26; K = "02/12/2023 02:13"; L = 1.75; M = "02/12/2023 12:59";
       N = 2.69; O = "02/12/2023 02:13"; P = 3.17; Q = "02/12/2023 12:59";
       R = 3.35; S = "02/12/2023 02:13"; T = 4.57; U = "02/12/2023 12:59";
       V = "https://www.betexplorer.com/football/serbia/prva-liga/metalac-kolubara/bLFrPzHH/" },

    @{ Row = 142; A = 141; E = 45262.54166666666; F = "Radnicki Beograd";  G = 1; H = "OFK Beograd";            I = 2;
       J = 5.34; K = "02/12/2023 01:13"; L = 5.25; M = "02/12/2023 12:56";
       N = 3.7;  O = "02/12/2023 01:13"; P = 3.88; Q = "02/12/2023 12:56";
       R = 1.53; S = "02/12/2023 01:13"; T = 1.52; U = "02/12/2023 12:56";
       V = "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-beograd-ofk-beograd/pSQLCgBl/" },

    @{ Row = 143; A = 142; E = 45262.54166666666; F = "Jedinstvo U.";      G = 1; H = "Tekstilac Odzaci";       I = 1;
       J = 2.33; K = "02/12/2023 01:13"; L = 2.19; M = "02/12/2023 12:19";
       N = 2.79; O = "02/12/2023 01:13"; P = 2.84; Q = "02/12/2023 12:19";
       R = 3.08; S = "02/12/2023 01:13"; T = 3.39; U = "02/12/2023 12:19";
       V = "https://www.betexplorer.com/football/serbia/prva-liga/jedinstvo-ub-tekstilac-odzaci/23LQBDQf/" },

    @{ Row = 144; A = 143; E = 45262.54166666666; F = "Sloboda";           G = 0; H = "Radnicki S. Mitrovica";  I = 1;
       J = 2.77; K = "02/12/2023 01:13"; L = 2.19; M = "02/12/2023 12:59";
       N = 2.69; O = "02/12/2023 01:13"; P = 2.89; Q = "02/12/2023 12:59";
       R = 2.63; S = "02/12/2023 01:13"; T = 3.31; U = "02/12/2023 12:59";
       V = "https://www.betexplorer.com/football/serbia/prva-liga/sloboda-radnicki-s-mitrovica/tfKUAXu1/" },

    @{ Row = 145; A = 144; E = 45262.625;         F = "Macva";             G = 2; H = "Smederevo";              I = 1;
       J = 2.29; K = "02/12/2023 03:12"; L = 2.01; M = "02/12/2023 13:53";
       N = 2.69; O = "02/12/2023 03:12"; P = 2.83; Q = "02/12/2023 14:40";
       R = 3.29; S = "02/12/2023 03:12"; T = 3.95; U = "02/12/2023 14:40";
       V = "https://www.betexplorer.com/football/serbia/prva-liga/macva-sabac-smederevo/lKNx9BAD/" },

    @{ Row = 146; A = 145; E = 45262.625;         F = "Dubocica";          G = 2; H = "Graficar Beograd";       I = 1;
       J = 2.33; K = "02/12/2023 03:12"; L = 2.5;  M = "02/12/2023 14:59";
       N = 2.97; O = "02/12/2023 03:12"; P = 3.03; Q = "02/12/2023 14:59";
       R = 2.87; S = "02/12/2023 03:12"; T = 2.67; U = "02/12/2023 14:59";
       V = "https://www.betexplorer.com/football/serbia/prva-liga/dubocica-graficar-beograd/6ySHDZer/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A$r").Value = $nr.A
    $ws.Range("B$r").Value = "serbia"
    $ws.Range("C$r").Value = "prva-liga"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $nr.E

    # F:V is 17 columns (F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T,U,V)
    $block = New-Object 'object[,]' 1,17
    $block[0,0]  = $nr.F
    $block[0,1]  = $nr.G
    $block[0,2]  = $nr.H
    $block[0,3]  = $nr.I
    $block[0,4]  = $nr.J
    $block[0,5]  = $nr.K
    $block[0,6]  = $nr.L
    $block[0,7]  = $nr.M
    $block[0,8]  = $nr.N
    $block[0,9]  = $nr.O
    $block[0,10] = $nr.P
    $block[0,11] = $nr.Q
    $block[0,12] = $nr.R
    $block[0,13] = $nr.S
    $block[0,14] = $nr.T
    $block[0,15] = $nr.U
    $block[0,16] = $nr.V
    Set-Block $ws $r $block
}

Write-Output "edit complete"
